# Add a "users" column to the "project hours" worksheet, listing the
# user(s) associated with each project row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New header cell, styled to match the existing header row (B1:D1)
$ws.Range("E1").Value = "users"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New data cells with the user list for each project
$ws.Range("E2").Value = "['Shuchen Song']"
$ws.Range("E3").Value = "['Hyung-Jin Yoon']"
